# Draft invoice generator: extend the hourly spot-price table from
# "25 aug 2022 - 31 aug 2022" so it covers the full week, and backfill
# the placeholder 0.00 values already present for 28 Aug (rows 74-97)
# with the real hourly prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill existing rows 74-97 (column B only; timestamps already there) ---
$existingValues = @(
    235.78, 46.24, 14.94, 14.34, 13.48, 12.61, 10.56, 14.19, 15.89, 17.46, 247.78, 266.5, 120.33, 78.64, 13.79, 64.49, 109.89, 255.06, 281.79, 242.82, 255.11, 253.58, 10.34, 8.91
)
for ($i = 0; $i -lt $existingValues.Length; $i++) {
    $row = 74 + $i
    $ws.Cells.Item($row, 2).Value = $existingValues[$i]
}

# --- Append new rows 98-169 for 29-31 Aug 2022: timestamp + price ---
$newRows = @(
    @("2022-08-29 00:00", 8.8),
    @("2022-08-29 01:00", 8.16),
    @("2022-08-29 02:00", 7.39),
    @("2022-08-29 03:00", 7.4),
    @("2022-08-29 04:00", 10.62),
    @("2022-08-29 05:00", 14.4),
    @("2022-08-29 06:00", 19.11),
    @("2022-08-29 07:00", 415.76),
    @("2022-08-29 08:00", 595.24),
    @("2022-08-29 09:00", 563.36),
    @("2022-08-29 10:00", 563.34),
    @("2022-08-29 11:00", 595.77),
    @("2022-08-29 12:00", 530.49),
    @("2022-08-29 13:00", 496.12),
    @("2022-08-29 14:00", 512.04),
    @("2022-08-29 15:00", 478.28),
    @("2022-08-29 16:00", 426.34),
    @("2022-08-29 17:00", 478.29),
    @("2022-08-29 18:00", 513.75),
    @("2022-08-29 19:00", 468.98),
    @("2022-08-29 20:00", 416.92),
    @("2022-08-29 21:00", 318.88),
    @("2022-08-29 22:00", 154.12),
    @("2022-08-29 23:00", 19.87),
    @("2022-08-30 00:00", 22.37),
    @("2022-08-30 01:00", 20.91),
    @("2022-08-30 02:00", 21.36),
    @("2022-08-30 03:00", 21.72),
    @("2022-08-30 04:00", 26.67),
    @("2022-08-30 05:00", 59.8),
    @("2022-08-30 06:00", 324.97),
    @("2022-08-30 07:00", 808.13),
    @("2022-08-30 08:00", 821.89),
    @("2022-08-30 09:00", 789.33),
    @("2022-08-30 10:00", 707.47),
    @("2022-08-30 11:00", 617.51),
    @("2022-08-30 12:00", 605.89),
    @("2022-08-30 13:00", 608.67),
    @("2022-08-30 14:00", 589.54),
    @("2022-08-30 15:00", 618.62),
    @("2022-08-30 16:00", 610.73),
    @("2022-08-30 17:00", 739.56),
    @("2022-08-30 18:00", 824.39),
    @("2022-08-30 19:00", 851.33),
    @("2022-08-30 20:00", 786.97),
    @("2022-08-30 21:00", 603.96),
    @("2022-08-30 22:00", 362.64),
    @("2022-08-30 23:00", 210.19),
    @("2022-08-31 00:00", 25),
    @("2022-08-31 01:00", 26.1),
    @("2022-08-31 02:00", 24.99),
    @("2022-08-31 03:00", 24.35),
    @("2022-08-31 04:00", 24.72),
    @("2022-08-31 05:00", 143.64),
    @("2022-08-31 06:00", 295.44),
    @("2022-08-31 07:00", 660.69),
    @("2022-08-31 08:00", 749.92),
    @("2022-08-31 09:00", 686.11),
    @("2022-08-31 10:00", 674.87),
    @("2022-08-31 11:00", 687.15),
    @("2022-08-31 12:00", 589.22),
    @("2022-08-31 13:00", 563.8),
    @("2022-08-31 14:00", 529.59),
    @("2022-08-31 15:00", 479.51),
    @("2022-08-31 16:00", 499.42),
    @("2022-08-31 17:00", 639.31),
    @("2022-08-31 18:00", 660.6),
    @("2022-08-31 19:00", 639.36),
    @("2022-08-31 20:00", 635.3),
    @("2022-08-31 21:00", 470.94),
    @("2022-08-31 22:00", 207.81),
    @("2022-08-31 23:00", 26.58)
)
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = 98 + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
}

# --- Grow the Table1 ListObject so it covers the newly added rows ---
# (the table ref has always run one row past the last data row, even in
# the original file: A1:B98 vs data ending at row 97 - keep that same
# one-row overhang now that data ends at row 169)
$lastRow = 97 + $newRows.Length
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:B" + ($lastRow + 1)))
